$wb = $excel.ActiveWorkbook

# Map of sheet name -> row -> new F-column value
$updates = @{
    "展览" = @{
        2  = 6754
        3  = 14
        9  = 95
        13 = 410
        15 = 1614
        17 = 3395
        21 = 2022
        22 = 140
        27 = 4
        28 = 135
    }
    "全部类型" = @{
        2  = 6754
        3  = 14
        10 = 95
        14 = 410
        16 = 1614
        18 = 3395
        22 = 2022
        23 = 140
        28 = 4
        29 = 135
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
